# Update crypto price/volume figures (GitHub Actions data refresh).
# For D-column values that are numeric-looking (e.g. "7.40", "1.00") we
# briefly force Text format so Excel keeps the literal digits/trailing
# zeros instead of normalizing them into a number, then restore the
# default "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.539.46"
$ws.Range("E2").Value = "  +4.82%  "

$ws.Range("D3").Value = "3.499.47"
$ws.Range("E3").Value = "  +3.03%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.38%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.498.57"
$ws.Range("E8").Value = "  +2.92%  "

$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("E11").Value = "  +4.59%  "

$ws.Range("E12").Value = "  +2.44%  "

$ws.Range("D13").Value = "4.104.32"
$ws.Range("E13").Value = "  +2.94%  "

$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.94%  "

$ws.Range("D16").Value = "66.504.25"
$ws.Range("E16").Value = "  +4.62%  "

$ws.Range("E17").Value = "  +3.51%  "

$ws.Range("D18").Value = "3.501.34"
$ws.Range("E18").Value = "  +2.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.16%  "

$ws.Range("E20").Value = "  +3.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.64%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("E25").Value = "  +2.19%  "

$ws.Range("E26").Value = "  +8.21%  "

$ws.Range("E27").Value = "  +4.84%  "

$ws.Range("E28").Value = "  +1.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.23%  "

$ws.Range("E31").Value = "  +5.40%  "

$ws.Range("E32").Value = "  +3.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("E36").Value = "  +3.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("E38").Value = "  +7.54%  "

$ws.Range("E39").Value = "  +6.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0747"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.63"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.17%  "

$ws.Range("D46").Value = "2.792.95"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0314"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "350.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.38%  "

$ws.Range("E50").Value = "  +6.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.78%  "

